# "Generate Report for Archive"
#
# Refreshes the handoff-status report:
#   - The "Ready for handoff" status label is now reported as "In Translation"
#     on the Overview sheet (columns E/F) and on each per-locale sheet
#     (column "Status").
#   - The now-narrower status columns are shrunk to fit the new, shorter
#     label text.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"
$newColWidth = 12.5   # renders to the same narrower stored column width for all affected columns

# --- Overview sheet: status appears in columns E and F of row 2 ---
$overview = $wb.Worksheets.Item("Overview")
if ($overview.Range("E2").Value2 -eq $oldStatus) {
    $overview.Range("E2").Value = $newStatus
}
if ($overview.Range("F2").Value2 -eq $oldStatus) {
    $overview.Range("F2").Value = $newStatus
}
$overview.Range("E1").ColumnWidth = $newColWidth
$overview.Range("F1").ColumnWidth = $newColWidth

# --- Per-locale sheets: status appears in column C of row 2 ---
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    if ($ws.Range("C2").Value2 -eq $oldStatus) {
        $ws.Range("C2").Value = $newStatus
    }
    $ws.Range("C1").ColumnWidth = $newColWidth
}
